{"js": "// Remove the trailing \"Ver no Jupiter ...\" / copyright footer block (and the\n// blank paragraph that separates it from the requisites list) that the\n// Jekyll site build no longer emits.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst VER_TEXT = \"Ver no Jupiter Salvar em pdf Salvar em docx\";\nconst COPYRIGHT_TEXT =\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\";\n\nconst items = paragraphs.items;\n\nlet verIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === VER_TEXT) {\n    verIndex = i;\n    break;\n  }\n}\n\nif (verIndex === -1) {\n  throw new Error('Could not locate paragraph with text \"' + VER_TEXT + '\"');\n}\n\nconst copyrightIndex = verIndex + 1;\nif (\n  copyrightIndex >= items.length ||\n  items[copyrightIndex].text !== COPYRIGHT_TEXT\n) {\n  throw new Error(\n    'Expected the copyright paragraph immediately after \"' + VER_TEXT + '\"'\n  );\n}\n\n// The blank paragraph right before \"Ver no Jupiter ...\" belongs to the same\n// footer block and should go too.\nlet startIndex = verIndex;\nif (verIndex - 1 >= 0 && items[verIndex - 1].text === \"\") {\n  startIndex = verIndex - 1;\n}\n\nconst paragraphsToDelete = [];\nfor (let i = startIndex; i <= copyrightIndex; i++) {\n  paragraphsToDelete.push(items[i]);\n}\n\n// Delete from the bottom up so earlier deletions don't disturb later ones.\nfor (let i = paragraphsToDelete.length - 1; i >= 0; i--) {\n  paragraphsToDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / copyright footer block (and the\n# blank paragraph that separates it from the requisites list) that the\n# Jekyll site build no longer emits.\n$d = $word.ActiveDocument\n\n$verText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$count = $d.Paragraphs.Count\n$verIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    $t = $t.TrimEnd([char]13, [char]7)\n    if ($t -eq $verText) {\n        $verIndex = $i\n        break\n    }\n}\n\nif ($verIndex -lt 1) {\n    throw \"Could not find paragraph with text: $verText\"\n}\n\n$copyIndex = $verIndex + 1\n$copyT = $d.Paragraphs.Item($copyIndex).Range.Text\n$copyT = $copyT.TrimEnd([char]13, [char]7)\nif ($copyT -ne $copyrightText) {\n    throw \"Expected the copyright paragraph immediately after '$verText'\"\n}\n\n# The blank paragraph right before \"Ver no Jupiter ...\" belongs to the same\n# footer block and should go too.\n$startIndex = $verIndex\nif ($verIndex -gt 1) {\n    $prevT = $d.Paragraphs.Item($verIndex - 1).Range.Text\n    $prevT = $prevT.TrimEnd([char]13, [char]7)\n    if ($prevT -eq \"\") {\n        $startIndex = $verIndex - 1\n    }\n}\n\n$startPara = $d.Paragraphs.Item($startIndex)\n$endPara = $d.Paragraphs.Item($copyIndex)\n$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)\n$rng.Delete()\n"}
